# Applies the numeric-value updates to the leve-profit tables across all
# eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), per the scheduled
# price-refresh run. Each sheet is backed by an Excel Table (Table_<JOB>)
# over A1:N141; columns H:N hold price/profit figures recomputed from the
# latest market data. A few rows also lose their HQ-profit figure (column N)
# entirely where the leve no longer has an HQ variant.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 461.3846
$ws.Range("I5").Value = 256.25
$ws.Range("J5").Value = 789.6
$ws.Range("K5").Value = 256.25
$ws.Range("L5").Value = 789.6
$ws.Range("M5").Value = -141.25
$ws.Range("N5").Value = -1019.6
$ws.Range("H17").Value = 2179555.5
$ws.Range("J17").Value = 2235280
$ws.Range("L17").Value = 6705840
$ws.Range("N17").Value = -6706176
$ws.Range("H87").Value = 48879.312
$ws.Range("J87").Value = 48879.312
$ws.Range("L87").Value = 48879.312
$ws.Range("N87").Value = -51375.312
$ws.Range("H90").Value = 48879.312
$ws.Range("J90").Value = 48879.312
$ws.Range("L90").Value = 146637.936
$ws.Range("N90").Value = -159117.936
$ws.Range("H135").Value = 1508.4615
$ws.Range("I135").Value = 1328.1818
$ws.Range("K135").Value = 11953.6362
$ws.Range("M135").Value = -9418.6362
$ws.Range("H138").Value = 1816.9375
$ws.Range("I138").Value = 1489.75
$ws.Range("J138").Value = 2798.5
$ws.Range("K138").Value = 4469.25
$ws.Range("L138").Value = 8395.5
$ws.Range("M138").Value = 670.75
$ws.Range("N138").Value = -18675.5
$ws.Range("H141").Value = 1234.25
$ws.Range("I141").Value = 1234.25
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3702.75
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 1477.25
$ws.Range("N141").ClearContents()

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 41920.668
$ws.Range("I28").Value = 32493
$ws.Range("K28").Value = 32493
$ws.Range("M28").Value = -32301
$ws.Range("H44").Value = 31657
$ws.Range("J44").Value = 31657
$ws.Range("L44").Value = 31657
$ws.Range("N44").Value = -32633
$ws.Range("H61").Value = 8861
$ws.Range("I61").Value = 6266.1665
$ws.Range("K61").Value = 6266.1665
$ws.Range("M61").Value = -6054.1665
$ws.Range("H97").Value = 1026.1
$ws.Range("I97").Value = 1175.826
$ws.Range("K97").Value = 1175.826
$ws.Range("M97").Value = -679.826
$ws.Range("H99").Value = 41920.668
$ws.Range("I99").Value = 32493
$ws.Range("K99").Value = 32493
$ws.Range("M99").Value = -29498
$ws.Range("H136").Value = 8861
$ws.Range("I136").Value = 6266.1665
$ws.Range("K136").Value = 18798.4995
$ws.Range("M136").Value = -16248.4995

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 111.9
$ws.Range("I22").Value = 115.71429
$ws.Range("K22").Value = 115.71429
$ws.Range("M22").Value = 57.28570999999999
$ws.Range("H35").Value = 96661
$ws.Range("J35").Value = 96661
$ws.Range("L35").Value = 96661
$ws.Range("N35").Value = -97281
$ws.Range("H107").Value = 616.5517
$ws.Range("I107").Value = 603.73914
$ws.Range("K107").Value = 603.73914
$ws.Range("M107").Value = 1316.26086
$ws.Range("H134").Value = 11666.167
$ws.Range("I134").Value = 6842
$ws.Range("K134").Value = 20526
$ws.Range("M134").Value = -17991

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 193.78125
$ws.Range("I22").Value = 191.46153
$ws.Range("K22").Value = 191.46153
$ws.Range("M22").Value = 158.53847
$ws.Range("H31").Value = 2836.56
$ws.Range("J31").Value = 4245.3335
$ws.Range("L31").Value = 4245.3335
$ws.Range("N31").Value = -4835.3335
$ws.Range("H34").Value = 2836.56
$ws.Range("J34").Value = 4245.3335
$ws.Range("L34").Value = 4245.3335
$ws.Range("N34").Value = -4649.3335
$ws.Range("H50").Value = 43567.855
$ws.Range("J50").Value = 43567.855
$ws.Range("L50").Value = 43567.855
$ws.Range("N50").Value = -44817.855
$ws.Range("H62").Value = 43574.68
$ws.Range("I62").Value = 72867.79
$ws.Range("J62").Value = 6292.5454
$ws.Range("K62").Value = 72867.79
$ws.Range("L62").Value = 6292.5454
$ws.Range("M62").Value = -72243.79
$ws.Range("N62").Value = -7540.5454
$ws.Range("H65").Value = 43574.68
$ws.Range("I65").Value = 72867.79
$ws.Range("J65").Value = 6292.5454
$ws.Range("K65").Value = 364338.95
$ws.Range("L65").Value = 31462.727
$ws.Range("M65").Value = -361218.95
$ws.Range("N65").Value = -37702.727
$ws.Range("H94").Value = 5719
$ws.Range("I94").Value = 5273.375
$ws.Range("J94").Value = 6164.625
$ws.Range("K94").Value = 5273.375
$ws.Range("L94").Value = 6164.625
$ws.Range("M94").Value = -4822.375
$ws.Range("N94").Value = -7066.625
$ws.Range("H122").Value = 1318.6522
$ws.Range("I122").Value = 1262.8182
$ws.Range("J122").Value = 1369.8334
$ws.Range("K122").Value = 3788.4546
$ws.Range("L122").Value = 4109.5002
$ws.Range("M122").Value = -1338.4546
$ws.Range("N122").Value = -9009.5002

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1372.12
$ws.Range("I5").Value = 615.36365
$ws.Range("J5").Value = 1966.7142
$ws.Range("K5").Value = 1846.09095
$ws.Range("L5").Value = 5900.142599999999
$ws.Range("M5").Value = -1734.09095
$ws.Range("N5").Value = -6124.142599999999
$ws.Range("H33").Value = 83.333336
$ws.Range("I33").Value = 75
$ws.Range("K33").Value = 450
$ws.Range("M33").Value = -167
$ws.Range("H128").Value = 500000
$ws.Range("I128").Value = 500000
$ws.Range("K128").Value = 1500000
$ws.Range("M128").Value = -1495020
$ws.Range("H135").Value = 1372.12
$ws.Range("I135").Value = 615.36365
$ws.Range("J135").Value = 1966.7142
$ws.Range("K135").Value = 5538.27285
$ws.Range("L135").Value = 17700.4278
$ws.Range("M135").Value = -3003.27285
$ws.Range("N135").Value = -22770.4278

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 29668.934
$ws.Range("J57").Value = 29668.934
$ws.Range("L57").Value = 29668.934
$ws.Range("N57").Value = -31308.934
$ws.Range("H102").Value = 1044.125
$ws.Range("I102").Value = 972.4706
$ws.Range("J102").Value = 1218.1428
$ws.Range("K102").Value = 972.4706
$ws.Range("L102").Value = 1218.1428
$ws.Range("M102").Value = 649.5294
$ws.Range("N102").Value = -4462.1428
$ws.Range("H132").Value = 13148.286
$ws.Range("I132").Value = 12571.869
$ws.Range("K132").Value = 37715.607
$ws.Range("M132").Value = -35185.607

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7600
$ws.Range("I122").Value = 7600
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 22800
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -20350
$ws.Range("N122").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5483.2173
$ws.Range("I122").Value = 2483.9412
$ws.Range("K122").Value = 7451.823600000001
$ws.Range("M122").Value = -5001.823600000001
$ws.Range("H132").Value = 17066.46
$ws.Range("I132").Value = 9544.355
$ws.Range("K132").Value = 28633.065
$ws.Range("M132").Value = -26103.065
$ws.Range("H136").Value = 2437.5652
$ws.Range("I136").Value = 2501.7083
$ws.Range("J136").Value = 2367.5908
$ws.Range("K136").Value = 7505.124899999999
$ws.Range("L136").Value = 7102.7724
$ws.Range("M136").Value = -4955.124899999999
$ws.Range("N136").Value = -12202.7724
